$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format before assignment so numeric-looking strings
# (e.g. "0.999", "1.00", "170.00") are preserved exactly as text and are
# not coerced into numeric values by Excel's type inference.
$cells = [ordered]@{
    "D2" = "65.080.89"
    "D3" = "3.478.39"
    "E3" = "  -0.79%  "
    "D4" = "0.999"
    "E4" = "  -0.04%  "
    "D5" = "588.11"
    "E5" = "  -2.66%  "
    "D6" = "137.04"
    "E6" = "  -4.29%  "
    "D7" = "3.477.41"
    "E7" = "  -0.76%  "
    "E8" = "  +0.05%  "
    "D9" = "0.490"
    "E9" = "  -2.75%  "
    "E10" = "  -5.32%  "
    "D11" = "7.13"
    "E11" = "  -7.41%  "
    "E12" = "  -4.54%  "
    "D13" = "4.066.11"
    "E13" = "  -0.74%  "
    "D14" = "0.0000180"
    "E14" = "  -6.54%  "
    "D15" = "3.479.65"
    "E15" = "  -0.73%  "
    "D16" = "26.58"
    "E16" = "  -6.71%  "
    "E17" = "  -1.26%  "
    "D18" = "65.048.55"
    "E18" = "  -1.91%  "
    "D19" = "9.72"
    "E19" = "  -8.13%  "
    "E20" = "  -5.06%  "
    "D21" = "13.93"
    "E21" = "  -4.19%  "
    "D22" = "388.41"
    "E22" = "  -7.46%  "
    "E23" = "  -4.77%  "
    "B24" = "LEO"
    "C24" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "D24" = "5.78"
    "E24" = "  +1.15%  "
    "B25" = "Dai"
    "C25" = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
    "D25" = "0.999"
    "E25" = "  -0.03%  "
    "B26" = "Litecoin"
    "C26" = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
    "D26" = "72.47"
    "E26" = "  -5.40%  "
    "D27" = "3.616.10"
    "E27" = "  -1.02%  "
    "E28" = "  -1.77%  "
    "E29" = "  -0.12%  "
    "D30" = "7.39"
    "E30" = "  -4.27%  "
    "D31" = "8.11"
    "E31" = "  -8.68%  "
    "D32" = "2.22"
    "E32" = "  -9.29%  "
    "D33" = "3.493.47"
    "E35" = "  -6.26%  "
    "D36" = "23.04"
    "E36" = "  -4.35%  "
    "D37" = "170.00"
    "E37" = "  -2.03%  "
    "B38" = "Aptos"
    "C38" = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
    "D38" = "6.83"
    "E38" = "  -8.12%  "
    "B39" = "Fetch.AI"
    "C39" = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
    "D39" = "1.19"
    "E39" = "  -9.44%  "
    "E40" = "  -9.32%  "
    "E41" = "  -8.40%  "
    "E42" = "  -3.11%  "
    "D43" = "0.811"
    "E43" = "  -4.37%  "
    "D44" = "42.57"
    "E44" = "  -6.37%  "
    "E45" = "  +0.08%  "
    "D46" = "24.91"
    "E46" = "  +9.30%  "
    "E47" = "  -11.79%  "
    "E48" = "  +5.29%  "
    "D49" = "1.61"
    "E49" = "  -7.96%  "
    "E50" = "  -4.29%  "
    "B51" = "dogwifhat"
    "C51" = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
    "D51" = "2.06"
    "E51" = "  -11.09%  "
}

foreach ($ref in $cells.Keys) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $cells[$ref]
}
